# Insert a new "Address" column before the existing "District" column (F),
# shifting District from F to G, and fill the new column with the school
# name + place portion of each record (everything from column B except the
# trailing district).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("F").Insert()

$ws.Range("F2").Value = 'Address'
$ws.Range("F3").Value = 'Sarvodya High SchoolShorapur'
$ws.Range("F4").Value = 'Dr. B S A AidedHigh SchoolShorapur'
$ws.Range("F5").Value = 'Govt. High School RajanakolurShorapur'
$ws.Range("F6").Value = 'Sagaranadu H S PetammapurShorapur'
$ws.Range("F7").Value = 'G H S DoranahalliShahapur'
$ws.Range("F8").Value = 'G H S Konkal'
$ws.Range("F9").Value = 'G H P School KurekuppaSandur'
$ws.Range("F10").Value = 'G H S talur siruguppa'
$ws.Range("F11").Value = 'Sri Agasara Dyavappa Govt High School Karur'
$ws.Range("F12").Value = 'G C P U CollegeEmmiganuree'
$ws.Range("F13").Value = 'S V G P U College HolaluHadagali'
$ws.Range("F14").Value = 'G M H P School KogaliHagari Bommanahalli'
$ws.Range("F15").Value = 'Govt. High School SirigeriSiruguppa'
$ws.Range("F16").Value = 'G H P S PoojarahalliKudligi'
$ws.Range("F17").Value = 'G H S GudekoteKudligi'
$ws.Range("F18").Value = 'G M H P S Vinobha Bhave ChithwadiHospet'
$ws.Range("F19").Value = 'G H S – Radio park'
$ws.Range("F20").Value = 'G H P S G N Hally'
$ws.Range("F21").Value = 'G H P S KarchiganurSiruguppa'
$ws.Range("F22").Value = 'Adarsha Vidyalaya (RMSA)Hospet'
$ws.Range("F23").Value = 'Govt. P U College High School sectionMariyamanahalliHospet'
$ws.Range("F24").Value = 'Govt. High School RaraviSiruguppa'
$ws.Range("F25").Value = 'Govt. High SchoolRavihalSiraguppa'
$ws.Range("F26").Value = 'MPM G H SAdavimallanakeriHoovinahadagali'
